$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $value) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.NumberFormat = "General"
    $c.Style = "Normal"
}

Set-TextValue 2 4 '26.242.26'
Set-TextValue 2 5 '  -4.12%  '
Set-TextValue 3 4 '1.656.57'
Set-TextValue 3 5 '  -3.58%  '
Set-TextValue 4 5 '  -0.18%  '
Set-TextValue 5 4 '216.32'
Set-TextValue 5 5 '  -3.84%  '
Set-TextValue 6 4 '0.5136'
Set-TextValue 6 5 '  -3.21%  '
Set-TextValue 7 5 '  -0.09%  '
Set-TextValue 8 4 '0.2601'
Set-TextValue 8 5 '  -2.27%  '
Set-TextValue 9 4 '0.06459'
Set-TextValue 9 5 '  -3.61%  '
Set-TextValue 10 4 '20.02'
Set-TextValue 10 5 '  -4.41%  '
Set-TextValue 11 4 '0.07796'
Set-TextValue 11 5 '  +1.23%  '
Set-TextValue 12 4 '1.660.89'
Set-TextValue 12 5 '  -3.28%  '
Set-TextValue 13 4 '4.301'
Set-TextValue 13 5 '  -4.12%  '
Set-TextValue 14 4 '1.884.22'
Set-TextValue 14 5 '  -3.59%  '
Set-TextValue 15 4 '0.5549'
Set-TextValue 15 5 '  -4.85%  '
Set-TextValue 16 4 '0.0₅8036'
Set-TextValue 16 5 '  -2.02%  '
Set-TextValue 17 4 '64.36'
Set-TextValue 17 5 '  -5.60%  '
Set-TextValue 18 4 '26.231.97'
Set-TextValue 18 5 '  -4.26%  '
Set-TextValue 19 2 'BitcoinCash'
Set-TextValue 19 3 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextValue 19 4 '211.27'
Set-TextValue 19 5 '  -5.10%  '
Set-TextValue 20 2 'Dai'
Set-TextValue 20 3 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue 20 4 '1.005'
Set-TextValue 20 5 '  -0.03%  '
Set-TextValue 21 4 '4.424'
Set-TextValue 21 5 '  -5.15%  '
Set-TextValue 22 4 '10.10'
Set-TextValue 22 5 '  -3.51%  '
Set-TextValue 23 4 '6.045'
Set-TextValue 23 5 '  +0.05%  '
Set-TextValue 24 5 '  -0.13%  '
Set-TextValue 25 4 '144.33'
Set-TextValue 25 5 '  -0.17%  '
Set-TextValue 26 4 '1.761'
Set-TextValue 26 5 '  +3.43%  '
Set-TextValue 27 5 '  -2.71%  '
Set-TextValue 28 4 '7.004'
Set-TextValue 28 5 '  -3.35%  '
Set-TextValue 29 4 '15.89'
Set-TextValue 29 5 '  -2.41%  '
Set-TextValue 30 4 '0.05113'
Set-TextValue 30 5 '  -5.05%  '
Set-TextValue 31 4 '1.246'
Set-TextValue 31 5 '  -3.87%  '
Set-TextValue 32 4 '3.368'
Set-TextValue 32 5 '  -3.26%  '
Set-TextValue 33 4 '3.226'
Set-TextValue 33 5 '  -5.90%  '
Set-TextValue 34 4 '1.564'
Set-TextValue 34 5 '  -4.14%  '
Set-TextValue 35 4 '2.742'
Set-TextValue 35 5 '  -4.31%  '
Set-TextValue 36 4 '0.9270'
Set-TextValue 36 5 '  -2.67%  '
Set-TextValue 37 4 '2.350'
Set-TextValue 37 5 '  -1.97%  '
Set-TextValue 38 4 '0.5733'
Set-TextValue 38 5 '  -2.70%  '
Set-TextValue 39 4 '1.170.05'
Set-TextValue 39 5 '  +1.32%  '
Set-TextValue 40 4 '0.01590'
Set-TextValue 40 5 '  -3.99%  '
Set-TextValue 41 4 '2.556'
Set-TextValue 41 5 '  -0.10%  '
Set-TextValue 42 5 '  -0.14%  '
Set-TextValue 43 4 '5.691'
Set-TextValue 43 5 '  -1.98%  '
Set-TextValue 44 4 '0.8249'
Set-TextValue 44 5 '  -2.07%  '
Set-TextValue 45 4 '100.41'
Set-TextValue 45 5 '  -0.65%  '
Set-TextValue 46 4 '1.794.88'
Set-TextValue 46 5 '  -3.56%  '
Set-TextValue 47 4 '0.0₈116'
Set-TextValue 47 5 '  +3.54%  '
Set-TextValue 48 4 '0.4547'
Set-TextValue 48 5 '  -0.96%  '
Set-TextValue 49 4 '55.53'
Set-TextValue 49 5 '  -4.13%  '
Set-TextValue 50 4 '1.006'
Set-TextValue 50 5 '  +0.20%  '
Set-TextValue 51 4 '7.879'
Set-TextValue 51 5 '  -3.47%  '
